# Updated cryptos list values per commit diff (price/volume refresh + row44/45 swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.156.62"
$ws.Range("E2").Value = "  +3.76%  "
$ws.Range("D3").Value = "'1.602.39"
$ws.Range("E3").Value = "  +3.13%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'213.03"
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  +2.38%  "
$ws.Range("E8").Value = "  +3.63%  "
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("D10").Value = "'18.02"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("E11").Value = "  +4.84%  "
$ws.Range("D12").Value = "'1.825.82"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("D13").Value = "'1.604.25"
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("D16").Value = "'26.146.48"
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("D17").Value = "'60.52"
$ws.Range("E17").Value = "  +3.50%  "
$ws.Range("D18").Value = "'0.0₃0721"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").Value = "'205.26"
$ws.Range("E20").Value = "  +11.65%  "
$ws.Range("E21").Value = "  +3.86%  "
$ws.Range("D22").Value = "'9.30"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("E23").Value = "  +3.15%  "
$ws.Range("E24").Value = "  +10.58%  "
$ws.Range("D25").Value = "'142.07"
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").Value = "'15.21"
$ws.Range("E28").Value = "  +3.53%  "
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("D31").Value = "'0.0471"
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("D32").Value = "'3.12"
$ws.Range("E32").Value = "  +4.22%  "
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("E34").Value = "  +2.55%  "
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("D36").Value = "'0.0163"
$ws.Range("E36").Value = "  +10.20%  "
$ws.Range("D37").Value = "'1.114.56"
$ws.Range("E37").Value = "  +2.92%  "
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("E40").Value = "  +3.75%  "
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("D42").Value = "'0.780"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("D43").Value = "'1.737.53"
$ws.Range("E43").Value = "  +3.21%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.12"
$ws.Range("E44").Value = "  +2.32%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'92.93"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  +5.70%  "
$ws.Range("D47").Value = "'53.39"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("D48").Value = "'0.0503"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").Value = "'0.408"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "'7.19"
$ws.Range("E51").Value = "  +1.50%  "
